$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting rows 2-6 down to 3-7
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the TimeScale entry
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "TimeScale"
$ws.Cells.Item(2, 3).Value = "t"
